$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a persistent run boundary at the edges of a Range by
# perturbing a formatting property and then restoring it. The runtime
# coalesces adjacent runs that resolve to identical formatting, but a
# boundary created this way survives even after the value is restored.
# ---------------------------------------------------------------------------
function New-RunBoundary($rng) {
    $rng.Font.Size = 20
    $rng.Font.Size = 16
}

# ---------------------------------------------------------------------------
# Edit 1: "...how many potential records we could have" -> add a trailing
# "." in its own run (same formatting).
# ---------------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("This was eventually changed to instead use the recordCount variable, since it holds how many potential records we could have", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos1 = $f1.End
$ins1 = $d.Range($endPos1, $endPos1)
$ins1.InsertAfter(".")
$newRun1 = $d.Range($endPos1, $endPos1 + 1)
New-RunBoundary $newRun1

# ---------------------------------------------------------------------------
# Edit 2: "...customer pointer at a time loaded, I decided..." ->
#         "...customer loaded at a time, I decided..."
# Produced as 5 runs:
#   "After realizing ... one customer "
#   "loaded "
#   "at a tim"
#   "e"
#   ", I decided to scrap ... before."
# ---------------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute("After realizing that the customer only wants to have one customer pointer at a time loaded, I decided to scrap the idea of texture caching, as it would still lead to the same problem as before.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $f2.Start

# delete "pointer " (8 chars right after "...one customer ")
$del2a = $d.Range($start2 + 66, $start2 + 66 + 8)
$del2a.Text = ""

# insert "loaded " in its place
$insPos2 = $start2 + 66
$ins2 = $d.Range($insPos2, $insPos2)
$ins2.InsertBefore("loaded ")
$run2b = $d.Range($insPos2, $insPos2 + 7)
New-RunBoundary $run2b

$afterLoadedPos2 = $insPos2 + 7   # start of "at a time loaded, I decided..."

# split "at a tim" | "e"
$eRng2 = $d.Range($afterLoadedPos2 + 8, $afterLoadedPos2 + 9)
New-RunBoundary $eRng2

# delete the stray " loaded" immediately after "e"
$del2b = $d.Range($afterLoadedPos2 + 9, $afterLoadedPos2 + 16)
$del2b.Text = ""

# give the tail (", I decided ... before.") its own run
$f2b = $d.Content
$f2b.Find.Execute(", I decided to scrap the idea of texture caching, as it would still lead to the same problem as before.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
New-RunBoundary $f2b

# ---------------------------------------------------------------------------
# Edit 3: "...reference a single Record pointer variable." ->
#         "...reference a single Record variable."
# ---------------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Record pointer variable.", $true, $false, $false, $false, $false, $true, 1, $false, "Record variable.", 2)

# ---------------------------------------------------------------------------
# Edit 4: "...I will need to  adapt the old load code..." (double space,
# 3 runs: "...need to ", " ", "adapt...") ->
#         "...I will need to adapt the old load code..." (single space,
# 3 runs: "...need ", "to adapt", " the old...")
# ---------------------------------------------------------------------------
$f4 = $d.Content
$f4.Find.Execute("need to  adapt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start4 = $f4.Start
$del4 = $d.Range($start4 + 7, $start4 + 8)
$del4.Text = ""

$mid4 = $d.Content
$mid4.Find.Execute("to adapt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
New-RunBoundary $mid4

# ---------------------------------------------------------------------------
# Edit 5: insert a new explanatory paragraph about GetRecord right after the
# two screenshots in the "Test 3" section (currently an empty paragraph).
# Formatting (sz=32/szCs=32/lang=en-US) is cloned from another run with
# identical rPr via Duplicate+FormattedText (copy, not move), then the
# cloned text is swapped for the real sentence with an in-place Find/Replace
# so the donor paragraph is left completely untouched.
# ---------------------------------------------------------------------------
$donor = $d.Content
$donor.Find.Execute("Add player was extremely simple to fix.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$donorDup = $donor.Duplicate
$template = $donorDup.FormattedText

$targetPara = $d.Paragraphs.Item(34)
$targetRng = $targetPara.Range
$gapPos = $targetRng.End - 1
$destRng = $d.Range($gapPos, $gapPos)
$destRng.FormattedText = $template

$targetPara2 = $d.Paragraphs.Item(34).Range
$targetPara2.Find.Execute("Add player was extremely simple to fix.", $true, $false, $false, $false, $false, $true, 1, $false, "Here the GetRecord function was reworked to search for the given index in the binary file. It would then set the Record to whatever it finds.", 2)

Write-Host "all edits applied"
